$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (shift existing A-D to B-E) ---
$ws.Columns("A").Insert()

# --- Insert a new row 2 (shift existing row 2 down to row 3) ---
$ws.Rows("2").Insert()

# --- Row 1 headers ---
$ws.Cells.Item(1,1).Value2 = "TabName"
$ws.Cells.Item(1,2).Value2 = "query"
$ws.Cells.Item(1,3).Value2 = "StatQuery"
$ws.Cells.Item(1,4).Value2 = "dbExcel"
$ws.Cells.Item(1,5).Value2 = "WebExcel"

# --- Column A tab-name labels (set first so shared strings stay grouped) ---
$ws.Cells.Item(2,1).Value2 = "CasesTab"
$ws.Cells.Item(3,1).Value2 = "FilesTab"

# --- Column B queries ---
$ws.Cells.Item(2,2).Value2 = 'MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN [''31504139'', ''31765263''] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity'
$ws.Cells.Item(3,2).Value2 = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN [''31504139'', ''31765263''] 
WITH
    f, parent, c, a, ct,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`'

# --- Column C stat queries (same text reused on both rows) ---
$ws.Cells.Item(2,3).Value2 = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
    WHERE a.pubmed_id IN [''31504139'', ''31765263''] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'
$ws.Cells.Item(3,3).Value2 = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
    WHERE a.pubmed_id IN [''31504139'', ''31765263''] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'

# --- Column D / E file names ---
$ws.Cells.Item(2,4).Value2 = "TC03_Trials_Filter_PubmedID-ALL_Neo4jData.xlsx"
$ws.Cells.Item(2,5).Value2 = "TC03_Trials_Filter_PubmedID-ALL_WebData.xlsx"
$ws.Cells.Item(3,4).Value2 = "TC03_Trials_Filter_PubmedID-ALL_Neo4jData.xlsx"
$ws.Cells.Item(3,5).Value2 = "TC03_Trials_Filter_PubmedID-ALL_WebData.xlsx"

# --- Row 4 : extra formatted (empty) cell C4 ---
$ws.Cells.Item(4,3).WrapText = $true

# --- Wrap text styling for the long query cells ---
$ws.Cells.Item(2,2).WrapText = $true
$ws.Cells.Item(2,3).WrapText = $true
$ws.Cells.Item(3,2).WrapText = $true
$ws.Cells.Item(3,3).WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# --- Column widths (nearest values achievable through this runtime's rounding) ---
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.6

# --- View / selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12:C13").Select()
